$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price (column D) updates: force Text format first so Excel does not
# coerce these numeric-looking strings into Double values, which would
# silently drop meaningful trailing/precision zeros (e.g. "246.70" -> 246.7).
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "246.70"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "26.32"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.069"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.05602"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "6.500"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.8135"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.02816"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.09379"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.001510"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0005966"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.006145"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.554"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06958"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.03116"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.1300"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.745"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04683"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.001249"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.004613"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.00009598"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.03660"
$ws.Range("B41").Value = "BKEXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.1366"
$ws.Range("E41").Value = "40BKEXTokenBKKBestin24h"
$ws.Range("B42").Value = "KickToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.006189"
$ws.Range("E42").Value = "41KickTokenKICK"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002659"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.008319"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005290"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00000000749"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.002062"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00002098"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0001998"
